# Update column F (CE) values in Sheet1 per data refresh ("Ensure GA works with lib.c")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "F3" = -2.272631571907541
    "F4" = -2.341895646539313
    "F5" = -2.425417410405929
    "F6" = -2.523196863507389
    "F7" = -2.620976316608849
    "F8" = -2.718755769710309
    "F9" = -2.830792912046613
    "F10" = -2.928572365148073
    "F11" = -3.040609507484376
    "F12" = -3.15264664982068
    "F13" = -3.264683792156983
    "F14" = -3.390978623728129
    "F15" = -3.531531144534118
    "F17" = -2.483697070400148
    "F18" = -2.502290548095691
    "F19" = -2.521990023379478
    "F20" = -2.541689498663263
    "F21" = -2.56138897394705
    "F22" = -2.582194446819077
    "F23" = -2.601893922102864
    "F24" = -2.622699394974892
    "F25" = -2.644546024313989
    "F26" = -2.666457494774261
    "F27" = -2.688368965234531
    "F28" = -2.711386433283045
    "F29" = -2.7355098989198
    "F30" = -2.752543291221651
    "F31" = -2.771516526297742
    "F32" = -2.790489761373831
    "F33" = -2.811608458051281
    "F34" = -2.830376074300251
    "F35" = -2.847409466602103
    "F36" = -2.866588320505312
    "F37" = -2.892096162153841
    "F38" = -2.909129554455695
    "F39" = -2.933922317854498
    "F40" = -2.946360946358146
    "F41" = -2.970948090929828
    "F42" = -2.991047641909888
    "F43" = -3.014410021957055
    "F44" = -3.033588875860268
    "F45" = -3.050720716782663
    "F46" = -3.079089098282275
    "F47" = -3.100207794959723
    "F48" = -3.126636333685097
    "F49" = -3.155004715184715
    "F50" = -3.165601949534669
    "F51" = -3.189474015856621
    "F52" = -3.212943566110371
    "F53" = -3.241723185264232
    "F54" = -3.262430644287436
    "F55" = -3.289528480774862
    "F56" = -3.31416140583422
    "F57" = -3.339736322140403
    "F58" = -3.35792740430911
    "F59" = -3.383640864784762
    "F60" = -3.41027502233726
    "F61" = -3.430982481360465
    "F62" = -3.453835401985027
    "F63" = -3.480263940710405
    "F64" = -3.502097715637579
    "F65" = -3.532200321084309
    "F66" = -3.553319017761757
    "F67" = -3.576887016636049
    "F68" = -3.602394858284577
    "F69" = -3.62734554071669
    "F70" = -3.653774079442065
    "F71" = -3.680044699134132
    "F309" = -2.938119656455225
    "F310" = -3.042370100245309
    "F311" = -3.132346051002539
    "F312" = -3.223306810957435
    "F313" = -3.33207975565473
    "F314" = -3.454089882385437
    "F315" = -3.584337650531226
    "F316" = -3.712581085015
    "F317" = -3.822220632181137
    "F318" = -3.952066137077082
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
